# Auto-generated edit script: update crypto price/volume values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.453.51"
$ws.Range("E2").Value = "  +12.76%  "
$ws.Range("D3").Value = "1.819.53"
$ws.Range("E3").Value = "  +7.77%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'231.23"
$ws.Range("E5").Value = "  +4.63%  "
$ws.Range("E6").Value = "  +4.83%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  +2.81%  "
$ws.Range("D9").Value = "'45.39"
$ws.Range("E9").Value = "  +2.75%  "
$ws.Range("E10").Value = "  +6.88%  "
$ws.Range("D11").Value = "'0.0684"
$ws.Range("E11").Value = "  +9.42%  "
$ws.Range("E12").Value = "  +3.32%  "
$ws.Range("D13").Value = "2.081.42"
$ws.Range("E13").Value = "  +7.77%  "
$ws.Range("D14").Value = "1.810.70"
$ws.Range("E14").Value = "  +7.26%  "
$ws.Range("D15").Value = "'0.648"
$ws.Range("E15").Value = "  +4.10%  "
$ws.Range("D16").Value = "34.345.86"
$ws.Range("E16").Value = "  +12.43%  "
$ws.Range("D17").Value = "'10.33"
$ws.Range("E17").Value = "  -4.23%  "
$ws.Range("E18").Value = "  +8.37%  "
$ws.Range("D19").Value = "'70.68"
$ws.Range("E19").Value = "  +7.09%  "
$ws.Range("D20").Value = "'261.63"
$ws.Range("E20").Value = "  +4.81%  "
$ws.Range("D21").Value = "0.0₃0753"
$ws.Range("E21").Value = "  +4.73%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("D23").Value = "'10.54"
$ws.Range("E23").Value = "  +3.62%  "
$ws.Range("E24").Value = "  +2.34%  "
$ws.Range("E25").Value = "  -1.16%  "
$ws.Range("D26").Value = "'161.28"
$ws.Range("E26").Value = "  +2.22%  "
$ws.Range("D27").Value = "'16.85"
$ws.Range("E27").Value = "  +5.60%  "
$ws.Range("E28").Value = "  +4.86%  "
$ws.Range("D29").Value = "'7.17"
$ws.Range("E29").Value = "  +5.71%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("E32").Value = "  +3.41%  "
$ws.Range("E33").Value = "  +6.86%  "
$ws.Range("E34").Value = "  +8.65%  "
$ws.Range("D35").Value = "1.588.10"
$ws.Range("E35").Value = "  +5.04%  "
$ws.Range("E36").Value = "  +5.44%  "
$ws.Range("E37").Value = "  +2.72%  "
$ws.Range("D38").Value = "'85.94"
$ws.Range("E38").Value = "  +8.15%  "
$ws.Range("D39").Value = "'0.633"
$ws.Range("E39").Value = "  +8.06%  "
$ws.Range("E40").Value = "  +5.35%  "
$ws.Range("D41").Value = "'2.82"
$ws.Range("E41").Value = "  +3.22%  "
$ws.Range("E42").Value = "  +7.86%  "
$ws.Range("D43").Value = "'2.34"
$ws.Range("E43").Value = "  +1.00%  "
$ws.Range("E44").Value = "  +7.34%  "
$ws.Range("D45").Value = "'0.0522"
$ws.Range("E45").Value = "  +3.48%  "
$ws.Range("E46").Value = "  +4.32%  "
$ws.Range("D47").Value = "1.980.77"
$ws.Range("E47").Value = "  +8.36%  "
$ws.Range("D48").Value = "'53.80"
$ws.Range("E48").Value = "  +2.75%  "
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("D51").Value = "'11.36"
$ws.Range("E51").Value = "  +18.34%  "
